$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.29910278320312
$ws.Range("C3").Value = 17.38476753234863
$ws.Range("C4").Value = 17.39287376403809
$ws.Range("C5").Value = 17.18902587890625
$ws.Range("C6").Value = 17.07291603088379
